$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Total" column header (column X = 24th column)
$ws.Cells.Item(1, 24).Value = "Total"

# Per-row totals for the existing disease-category rows (2-6)
$ws.Cells.Item(2, 24).Value = 2079
$ws.Cells.Item(3, 24).Value = 302
$ws.Cells.Item(4, 24).Value = 1084
$ws.Cells.Item(5, 24).Value = 239
$ws.Cells.Item(6, 24).Value = 1385

# New row 7: "Outros" category
$ws.Cells.Item(7, 1).Value = "Outros"
$ws.Cells.Item(7, 2).Value = 136
$ws.Cells.Item(7, 3).Value = 2
$ws.Cells.Item(7, 4).Value = 10
$ws.Cells.Item(7, 5).Value = 28
$ws.Cells.Item(7, 6).Value = 80
$ws.Cells.Item(7, 7).Value = 71
$ws.Cells.Item(7, 8).Value = 75
$ws.Cells.Item(7, 9).Value = 94
$ws.Cells.Item(7, 10).Value = 80
$ws.Cells.Item(7, 11).Value = 82
$ws.Cells.Item(7, 12).Value = 118
$ws.Cells.Item(7, 13).Value = 111
$ws.Cells.Item(7, 14).Value = 114
$ws.Cells.Item(7, 15).Value = 123
$ws.Cells.Item(7, 16).Value = 122
$ws.Cells.Item(7, 17).Value = 154
$ws.Cells.Item(7, 18).Value = 197
$ws.Cells.Item(7, 19).Value = 194
$ws.Cells.Item(7, 20).Value = 124
$ws.Cells.Item(7, 21).Value = 45
$ws.Cells.Item(7, 22).Value = 7
$ws.Cells.Item(7, 23).Value = 1
$ws.Cells.Item(7, 24).Value = 1968

# New row 8: "Total" category (grand total row)
$ws.Cells.Item(8, 1).Value = "Total"
$ws.Cells.Item(8, 2).Value = 158
$ws.Cells.Item(8, 3).Value = 6
$ws.Cells.Item(8, 4).Value = 13
$ws.Cells.Item(8, 5).Value = 39
$ws.Cells.Item(8, 6).Value = 98
$ws.Cells.Item(8, 7).Value = 91
$ws.Cells.Item(8, 8).Value = 113
$ws.Cells.Item(8, 9).Value = 170
$ws.Cells.Item(8, 10).Value = 173
$ws.Cells.Item(8, 11).Value = 239
$ws.Cells.Item(8, 12).Value = 360
$ws.Cells.Item(8, 13).Value = 476
$ws.Cells.Item(8, 14).Value = 535
$ws.Cells.Item(8, 15).Value = 635
$ws.Cells.Item(8, 16).Value = 700
$ws.Cells.Item(8, 17).Value = 806
$ws.Cells.Item(8, 18).Value = 859
$ws.Cells.Item(8, 19).Value = 861
$ws.Cells.Item(8, 20).Value = 504
$ws.Cells.Item(8, 21).Value = 178
$ws.Cells.Item(8, 22).Value = 42
$ws.Cells.Item(8, 23).Value = 1
$ws.Cells.Item(8, 24).Value = 7057
